$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: the previously-blank C5/D5/G5/H5 cells used a redundant cell
# style (cellXfs index 17) that is identical in appearance to the style
# already used elsewhere (index 5 - centered text, border 3). Re-apply
# that formatting so the duplicate style collapses away.
# (Two calls instead of one multi-area "C5:D5,G5:H5" range, since only
# the first area of a multi-area range picks up the format.)
$ws.Range("C5:D5").HorizontalAlignment = -4108
$ws.Range("C5:D5").VerticalAlignment = -4108
$ws.Range("G5:H5").HorizontalAlignment = -4108
$ws.Range("G5:H5").VerticalAlignment = -4108

# Row 12 "Max D wire": D12/E12/F12 move from 0.33 to 0.5 (matches C12),
# which in turn changes the dependent "Max turns possible" formulas in
# row 13 (D13/E13/F13 recalc from 137.74... to 60).
$ws.Range("D12:F12").Value = 0.5

# Update the active selection to E5 (matches the saved workbook view).
$ws.Range("E5").Select() | Out-Null
